$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.01132367786385012, 2.289151444524298, 20.74398971997876, 4.554557027854494, 4.656905011860751, 23),
    @(-0.0664853343312105, 2.410518245823468, 18.500884502483, 4.301265453617459, 4.401959587781231, 22),
    @(-0.5622320716953093, 1.875594170520679, 10.40438817800559, 3.225583385684764, 3.25464244416202, 21),
    @(-0.1812427017839789, 1.441856625662851, 8.636573577577872, 2.938804787252442, 3.009410619098892, 20),
    @(-0.2071371915338822, 1.880800794843898, 13.16300375377707, 3.628085411587917, 3.721423419543999, 19),
    @(-0.2618932902003571, 1.772026561808276, 9.805127963275318, 3.131314095276186, 3.21080636400065, 18),
    @(-0.2336496272831486, 1.757688402756847, 11.70993972551334, 3.4219789195016, 3.519063341396083, 17),
    @(-0.110650871936715, 1.929024814611791, 10.40150795858265, 3.225136889898264, 3.328946074360039, 16),
    @(-0.1706632301364332, 2.15590519481267, 15.75989337962914, 3.969873219591419, 4.105410294639327, 15),
    @(-0.1306111199879883, 1.99991372347226, 10.84813572458541, 3.293650820075713, 3.415294451084287, 14)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
